$wb = $excel.ActiveWorkbook

# 1. "ms_source list" sheet: insert a new value "DESI" before the existing
#    "nanoDESI" entry (row 7), which pushes nanoDESI down to row 8.
$msSourceWs = $wb.Worksheets.Item("ms_source list")
$msSourceWs.Rows.Item(7).Insert()
$msSourceWs.Range("A7").Value = "DESI"

# 2. Update the data validation on the "Export as TSV" sheet (column Q,
#    ms_source) so its source range covers the new list size (A1:A8).
$mainWs = $wb.Worksheets.Item("Export as TSV")
$qValidationRange = $mainWs.Range("Q2:Q1048576")
$newFormula = '''ms_source list''!$A$1:$A$8'
$qValidationRange.Validation.Formula1 = $newFormula
